$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.39526366666666
$ws.Range("H2").Value = 280.185791
$ws.Range("I2").Value = 0.2167755775732346
$ws.Range("J2").Value = 0.2167755775732346
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 8546.722396692241
$ws.Range("R2").Value = 76920.50157023018
$ws.Range("S2").Value = 0.2099584153179167
$ws.Range("T2").Value = 0.2099584153179167
$ws.Range("G3").Value = 93.39526366666666
$ws.Range("H3").Value = 280.185791
$ws.Range("I3").Value = 0.2167755775732346
$ws.Range("J3").Value = 0.2167755775732346
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 12.581462759064
$ws.Range("R3").Value = 113.233164831576
$ws.Range("S3").Value = 0.0003090756737690238
$ws.Range("T3").Value = 0.0003090756737690239
$ws.Range("G4").Value = 93.39526366666666
$ws.Range("H4").Value = 280.185791
$ws.Range("I4").Value = 0.2167755775732346
$ws.Range("J4").Value = 0.2167755775732346
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 264.9229813528206
$ws.Range("R4").Value = 2384.306832175385
$ws.Range("S4").Value = 0.006508086581548896
$ws.Range("T4").Value = 0.006508086581548895
$ws.Range("I5").Value = 0.5566060939249745
$ws.Range("J5").Value = 0.5566060939249745
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 21945.0817400167
$ws.Range("R5").Value = 197505.7356601504
$ws.Range("S5").Value = 0.539101935490414
$ws.Range("T5").Value = 0.539101935490414
$ws.Range("I6").Value = 0.5566060939249745
$ws.Range("J6").Value = 0.5566060939249745
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.0007936014076386764
$ws.Range("T6").Value = 0.0007936014076386765
$ws.Range("I7").Value = 0.5566060939249745
$ws.Range("J7").Value = 0.5566060939249745
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("Q7").Value = 680.2322821256734
$ws.Range("R7").Value = 6122.09053913106
$ws.Range("S7").Value = 0.01671055702692191
$ws.Range("T7").Value = 0.01671055702692191
$ws.Range("I8").Value = 0.226618328501791
$ws.Range("J8").Value = 0.2266183285017909
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 8934.788528255131
$ws.Range("R8").Value = 80413.09675429617
$ws.Range("S8").Value = 0.2194916312385639
$ws.Range("T8").Value = 0.2194916312385639
$ws.Range("I9").Value = 0.226618328501791
$ws.Range("J9").Value = 0.2266183285017909
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("S9").Value = 0.0003231093343365133
$ws.Range("T9").Value = 0.0003231093343365133
$ws.Range("I10").Value = 0.226618328501791
$ws.Range("J10").Value = 0.2266183285017909
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("S10").Value = 0.006803587928890594
$ws.Range("T10").Value = 0.006803587928890593
